$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire row 11 (1998 - Lodoss Island World Guide), then delete it,
# shifting all rows below up by one.
$ws.Rows.Item(11).Select() | Out-Null
$ws.Rows.Item(11).Delete()
